$wb = $excel.ActiveWorkbook

# Update OFF sheet (Week totals through current week) - row 2 ("H")
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B2").Value = 193
$wsOff.Range("C2").Value = 118
$wsOff.Range("D2").Value = 52
$wsOff.Range("E2").Value = 25

# Update DEF sheet (Week totals through current week) - row 2 ("H")
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B2").Value = 400
$wsDef.Range("C2").Value = 289
$wsDef.Range("D2").Value = 78
$wsDef.Range("F2").Value = 6
$wsDef.Range("G2").Value = 3
